# Auto-generated edit script: updates recalculated market price/profit
# figures on several Leve rows across all 8 Sheets tabs.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3097.082
$ws.Range("I15").Value = 3097.082
$ws.Range("K15").Value = 9291.245999999999
$ws.Range("M15").Value = -9122.245999999999
$ws.Range("H41").Value = 299.75
$ws.Range("I41").Value = 215.6
$ws.Range("J41").Value = 440
$ws.Range("K41").Value = 215.6
$ws.Range("L41").Value = 440
$ws.Range("M41").Value = 224.4
$ws.Range("N41").Value = -1320
$ws.Range("H58").Value = 1068
$ws.Range("I58").Value = 750.05884
$ws.Range("J58").Value = 1668.5555
$ws.Range("K58").Value = 2250.17652
$ws.Range("L58").Value = 5005.666499999999
$ws.Range("M58").Value = -2100.17652
$ws.Range("N58").Value = -5305.666499999999
$ws.Range("H62").Value = 5333.4116
$ws.Range("I62").Value = 4125
$ws.Range("K62").Value = 4125
$ws.Range("M62").Value = -3501
$ws.Range("H63").Value = 29750
$ws.Range("J63").Value = 29750
$ws.Range("L63").Value = 29750
$ws.Range("N63").Value = -30998
$ws.Range("H65").Value = 5333.4116
$ws.Range("I65").Value = 4125
$ws.Range("K65").Value = 20625
$ws.Range("M65").Value = -17505
$ws.Range("H66").Value = 29750
$ws.Range("J66").Value = 29750
$ws.Range("L66").Value = 89250
$ws.Range("N66").Value = -95490
$ws.Range("H80").Value = 10172.5
$ws.Range("I80").Value = 8667.333000000001
$ws.Range("J80").Value = 10473.533
$ws.Range("K80").Value = 26001.999
$ws.Range("L80").Value = 31420.599
$ws.Range("M80").Value = -25003.999
$ws.Range("N80").Value = -33416.599
$ws.Range("H83").Value = 10172.5
$ws.Range("I83").Value = 8667.333000000001
$ws.Range("J83").Value = 10473.533
$ws.Range("K83").Value = 78005.997
$ws.Range("L83").Value = 94261.79699999999
$ws.Range("M83").Value = -73013.997
$ws.Range("N83").Value = -104245.797
$ws.Range("H98").Value = 914.7727
$ws.Range("I98").Value = 914.7727
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 914.7727
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 583.2273
$ws.Range("N98").ClearContents()
$ws.Range("H99").Value = 1382.9231
$ws.Range("I99").Value = 1081.5
$ws.Range("J99").Value = 5000
$ws.Range("K99").Value = 3244.5
$ws.Range("L99").Value = 15000
$ws.Range("M99").Value = -1746.5
$ws.Range("N99").Value = -17996
$ws.Range("H101").Value = 2129.9092
$ws.Range("I101").Value = 527.53845
$ws.Range("J101").Value = 4444.4443
$ws.Range("K101").Value = 1582.61535
$ws.Range("L101").Value = 13333.3329
$ws.Range("M101").Value = 39.38464999999997
$ws.Range("N101").Value = -16577.3329
$ws.Range("H122").Value = 914.7727
$ws.Range("I122").Value = 914.7727
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2744.3181
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -294.3181
$ws.Range("N122").ClearContents()
$ws.Range("H125").Value = 933.3333
$ws.Range("I125").Value = 850
$ws.Range("J125").Value = 975
$ws.Range("K125").Value = 7650
$ws.Range("L125").Value = 8775
$ws.Range("M125").Value = -5190
$ws.Range("N125").Value = -13695
$ws.Range("H132").Value = 3134.9211
$ws.Range("I132").Value = 1872.92
$ws.Range("J132").Value = 5561.846
$ws.Range("K132").Value = 5618.76
$ws.Range("L132").Value = 16685.538
$ws.Range("M132").Value = -3088.76
$ws.Range("N132").Value = -21745.538

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 5385.3076
$ws.Range("I30").Value = 1701.5
$ws.Range("J30").Value = 8542.857
$ws.Range("K30").Value = 1701.5
$ws.Range("L30").Value = 8542.857
$ws.Range("M30").Value = -1551.5
$ws.Range("N30").Value = -8842.857
$ws.Range("H32").Value = 4323.655
$ws.Range("I32").Value = 3770.0513
$ws.Range("J32").Value = 5673.0625
$ws.Range("K32").Value = 3770.0513
$ws.Range("L32").Value = 5673.0625
$ws.Range("M32").Value = -3483.0513
$ws.Range("N32").Value = -6247.0625
$ws.Range("H52").Value = 17999.5
$ws.Range("J52").Value = 17999.5
$ws.Range("L52").Value = 17999.5
$ws.Range("N52").Value = -18635.5
$ws.Range("H97").Value = 980.1177
$ws.Range("I97").Value = 938.5714
$ws.Range("J97").Value = 1174
$ws.Range("K97").Value = 938.5714
$ws.Range("L97").Value = 1174
$ws.Range("M97").Value = -442.5714
$ws.Range("N97").Value = -2166
$ws.Range("H119").Value = 29554.363
$ws.Range("J119").Value = 29554.363
$ws.Range("L119").Value = 29554.363
$ws.Range("N119").Value = -39230.363
$ws.Range("H121").Value = 30000
$ws.Range("J121").Value = 30000
$ws.Range("L121").Value = 30000
$ws.Range("N121").Value = -33494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 562.6923
$ws.Range("I22").Value = 502.33334
$ws.Range("J22").Value = 698.5
$ws.Range("K22").Value = 502.33334
$ws.Range("L22").Value = 698.5
$ws.Range("M22").Value = -152.33334
$ws.Range("N22").Value = -1398.5
$ws.Range("H23").Value = 18653
$ws.Range("I23").Value = 10000
$ws.Range("J23").Value = 22979.5
$ws.Range("K23").Value = 10000
$ws.Range("L23").Value = 22979.5
$ws.Range("M23").Value = -9760
$ws.Range("N23").Value = -23459.5
$ws.Range("H27").Value = 18653
$ws.Range("I27").Value = 10000
$ws.Range("J27").Value = 22979.5
$ws.Range("K27").Value = 10000
$ws.Range("L27").Value = 22979.5
$ws.Range("M27").Value = -9808
$ws.Range("N27").Value = -23363.5
$ws.Range("H62").Value = 2771.4285
$ws.Range("I62").Value = 2250
$ws.Range("J62").Value = 3710
$ws.Range("K62").Value = 2250
$ws.Range("L62").Value = 3710
$ws.Range("M62").Value = -1626
$ws.Range("N62").Value = -4958
$ws.Range("H65").Value = 2771.4285
$ws.Range("I65").Value = 2250
$ws.Range("J65").Value = 3710
$ws.Range("K65").Value = 11250
$ws.Range("L65").Value = 18550
$ws.Range("M65").Value = -8130
$ws.Range("N65").Value = -24790
$ws.Range("H134").Value = 1648
$ws.Range("I134").Value = 1739.7
$ws.Range("K134").Value = 5219.1
$ws.Range("M134").Value = -2684.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2590
$ws.Range("J39").Value = 2817.1428
$ws.Range("L39").Value = 8451.428400000001
$ws.Range("N39").Value = -9039.428400000001
$ws.Range("H131").Value = 971.9298
$ws.Range("J131").Value = 971.9298
$ws.Range("L131").Value = 2915.7894
$ws.Range("N131").Value = -12995.7894

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 21459.715
$ws.Range("I25").Value = 1550
$ws.Range("J25").Value = 48006
$ws.Range("K25").Value = 1550
$ws.Range("L25").Value = 48006
$ws.Range("M25").Value = -1021
$ws.Range("N25").Value = -49064
$ws.Range("H121").Value = 20000
$ws.Range("J121").Value = 20000
$ws.Range("L121").Value = 20000
$ws.Range("N121").Value = -23494
$ws.Range("H122").Value = 28573.676
$ws.Range("I122").Value = 48806.145
$ws.Range("J122").Value = 2018.5625
$ws.Range("K122").Value = 146418.435
$ws.Range("L122").Value = 6055.6875
$ws.Range("M122").Value = -143968.435
$ws.Range("N122").Value = -10955.6875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 887.9474
$ws.Range("I46").Value = 998.875
$ws.Range("J46").Value = 807.2727
$ws.Range("K46").Value = 998.875
$ws.Range("L46").Value = 807.2727
$ws.Range("M46").Value = -810.875
$ws.Range("N46").Value = -1183.2727
$ws.Range("H119").Value = 30105
$ws.Range("J119").Value = 30105
$ws.Range("L119").Value = 30105
$ws.Range("N119").Value = -39781

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 10000
$ws.Range("J30").Value = 10000
$ws.Range("L30").Value = 10000
$ws.Range("N30").Value = -10214
$ws.Range("H119").Value = 29528.285
$ws.Range("J119").Value = 29528.285
$ws.Range("L119").Value = 29528.285
$ws.Range("N119").Value = -39204.285
$ws.Range("H122").Value = 16748033
$ws.Range("I122").Value = 31250960
$ws.Range("J122").Value = 173258.58
$ws.Range("K122").Value = 93752880
$ws.Range("L122").Value = 519775.74
$ws.Range("M122").Value = -93750430
$ws.Range("N122").Value = -524675.74
$ws.Range("H132").Value = 3574523.8
$ws.Range("I132").Value = 4258296.5
$ws.Range("J132").Value = 3710.5557
$ws.Range("K132").Value = 12774889.5
$ws.Range("L132").Value = 11131.6671
$ws.Range("M132").Value = -12772359.5
$ws.Range("N132").Value = -16191.6671

